$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 2: "Prompts provided to the LLM for each choice:" write-up (rows 89-123) ---
# Cell text values, written in the same order the original author typed them so the
# shared-string table comes out in the same sequence as the authored workbook.
$ws.Range('C89').Value = 'Prompts provided to the LLM for each choice:'
$ws.Range('C91').Value = 'Flat Line:'
$ws.Range('C100').Value = 'Distribute {qty entered by user} items across {number of months entered by user} months '
$ws.Range('C92').Value = 'Distribute {qty entered by user} items across {number of months entered by user} months'
$ws.Range('C93').Value = 'using a flat line so that each month has either approximately or exactly the same number'
$ws.Range('C94').Value = 'of items as all other months. All deliveries must be whole numbers and the total needs'
$ws.Range('C95').Value = 'to add to {qty entered by user}. Do not include any text describing the result, respond '
$ws.Range('C101').Value = 'using an ascending line so that each month has the same or more items than the month before '
$ws.Range('C102').Value = 'it. Try to get the line as close to a 45 degree angle as possible. All monthly quantities '
$ws.Range('C103').Value = 'must be positive integers and the total needs to add to {qty entered by user}.  Do not '
$ws.Range('C104').Value = 'include any text describing the result, respond with only an array containing the quantity '
$ws.Range('C105').Value = 'for each month like this example [3, 5, 2].'
$ws.Range('C99').Value = 'Ramp Up:'
$ws.Range('C96').Value = 'with only an array containing the quantity for each month like this example [3, 5, 2].'
$ws.Range('C108').Value = 'Ramp Down:'
$ws.Range('C110').Value = 'using a descending line so that each month has the same or fewer items than the month before '
$ws.Range('C117').Value = 'Bell Curve:'
$ws.Range('C118').Value = 'Distribute {qty entered by user} items across {number of months entered by user} months using'
$ws.Range('C119').Value = 'a bell curve so that the month or months in the middle have the most items and the months on '
$ws.Range('C120').Value = 'the ends have the fewest items. Try to make the bell curve not too steep or too flat. All '
$ws.Range('C121').Value = 'monthly quantities must be whole numbers and the total needs to add to {qty entered by user}.  '
$ws.Range('C122').Value = 'Do not include any text describing the result, respond with only an array containing the '
$ws.Range('C123').Value = 'quantity for each month like this example [3, 5, 2].'
$ws.Range('C109').Value = 'Distribute {qty entered by user} items across {number of months entered by user} months '
$ws.Range('C111').Value = 'it. Try to get the line as close to a 45 degree angle as possible. All monthly quantities '
$ws.Range('C112').Value = 'must be positive integers and the total needs to add to {qty entered by user}.  Do not '
$ws.Range('C113').Value = 'include any text describing the result, respond with only an array containing the quantity '
$ws.Range('C114').Value = 'for each month like this example [3, 5, 2].'

# --- Formatting ---
# Row 92 carries a custom row height and was touched for alignment (matches the
# "applyAlignment" style seen with no explicit alignment override).
$ws.Range('C92:L92').WrapText = $false
$ws.Range('C92:L92').RowHeight = 13.8

# Bold section headers.
$ws.Range('C89:G89').Font.Bold = $true
$ws.Range('C91').Font.Bold = $true
$ws.Range('C99').Font.Bold = $true
$ws.Range('C108').Font.Bold = $true
$ws.Range('C117').Font.Bold = $true

# --- View state ---
$null = $ws.Range('K100').Select()
